$d = $word.ActiveDocument

# 1. "shall use JDBC to connect to a MariaDB" -> "shall connect to a MariaDB"
#    (drop " use JDBC to", keep the rest)
$r1 = $d.Content.Find.Execute("shall use JDBC to connect", $true, $false, $false, $false, $false, $true, 1, $false, "shall connect", 2)

# 2. Insert a new sentence about using Hibernate instead of JDBC, right between
#    "database. " and "The middle tier"
$r2 = $d.Content.Find.Execute("database. The middle tier", $true, $false, $false, $false, $false, $true, 1, $false, "database. Instead of utilizing JDBC, you should be utilizing Hibernate to perform database operations. The middle tier", 2)

# 3. Drop " for dynamic Web application development" (keep the trailing period).
#    The search text intentionally excludes "Javalin" so its spell-check
#    proofErr wrapper (<w:proofErr w:type="spellStart/spellEnd"/>) is left intact.
$r3 = $d.Content.Find.Execute(" technology for dynamic Web application development.", $true, $false, $false, $false, $false, $true, 1, $false, " technology.", 2)

# 4. Front-end description rewrite: vanilla JavaScript (multi-page) vs Angular (single page)
$r4 = $d.Content.Find.Execute("can use JavaScript or Angular to make a single page application that uses AJAX", $true, $false, $false, $false, $false, $true, 1, $false, "can use vanilla JavaScript (multi-page application) or Angular (single page application) that uses AJAX", 2)

# 5. Passwords: "encrypted in Java" -> "hashed"
$r5 = $d.Content.Find.Execute("Passwords shall be encrypted in Java and securely stored", $true, $false, $false, $false, $false, $true, 1, $false, "Passwords shall be hashed and securely stored", 2)

Write-Output "replace results: $r1 $r2 $r3 $r4 $r5"
